$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price cells that are about to be rewritten to be treated as
# text, so numeric-looking price strings (e.g. "301.65", "0.530", "1.00")
# keep their exact text representation instead of being auto-converted to
# numbers by Excel.
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D7").NumberFormat = "@"
$ws.Range("D9:D12").NumberFormat = "@"
$ws.Range("D14:D20").NumberFormat = "@"
$ws.Range("D22:D26").NumberFormat = "@"
$ws.Range("D28:D33").NumberFormat = "@"
$ws.Range("D35:D36").NumberFormat = "@"
$ws.Range("D39:D44").NumberFormat = "@"
$ws.Range("D47:D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "41.871.42"
$ws.Range("E2").Value = "  +5.42%  "

$ws.Range("D3").Value = "2.254.61"
$ws.Range("E3").Value = "  +1.66%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "301.65"
$ws.Range("E5").Value = "  +3.64%  "

$ws.Range("D6").Value = "92.18"
$ws.Range("E6").Value = "  +6.23%  "

$ws.Range("D7").Value = "0.530"
$ws.Range("E7").Value = "  +3.50%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "0.481"
$ws.Range("E9").Value = "  +3.09%  "

$ws.Range("D10").Value = "54.55"
$ws.Range("E10").Value = "  +9.36%  "

$ws.Range("D11").Value = "32.40"
$ws.Range("E11").Value = "  +6.35%  "

$ws.Range("D12").Value = "0.0797"
$ws.Range("E12").Value = "  +2.19%  "

$ws.Range("E13").Value = "  +2.96%  "

$ws.Range("D14").Value = "6.66"
$ws.Range("E14").Value = "  +3.23%  "

$ws.Range("D15").Value = "2.605.07"
$ws.Range("E15").Value = "  +1.86%  "

$ws.Range("D16").Value = "14.11"
$ws.Range("E16").Value = "  +2.47%  "

$ws.Range("D17").Value = "2.270.56"
$ws.Range("E17").Value = "  +1.82%  "

$ws.Range("D18").Value = "0.755"
$ws.Range("E18").Value = "  +3.52%  "

$ws.Range("D19").Value = "41.775.60"
$ws.Range("E19").Value = "  +5.26%  "

$ws.Range("D20").Value = "12.12"
$ws.Range("E20").Value = "  +9.37%  "

$ws.Range("E21").Value = "  +1.74%  "

$ws.Range("D22").Value = "5.92"
$ws.Range("E22").Value = "  +3.20%  "

$ws.Range("D23").Value = "67.03"
$ws.Range("E23").Value = "  +2.15%  "

$ws.Range("D24").Value = "240.60"
$ws.Range("E24").Value = "  +1.49%  "

$ws.Range("D25").Value = "2.56"
$ws.Range("E25").Value = "  +4.73%  "

$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("E27").Value = "  +3.87%  "

$ws.Range("D28").Value = "23.86"
$ws.Range("E28").Value = "  +3.65%  "

$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  +3.04%  "

$ws.Range("D30").Value = "9.64"
$ws.Range("E30").Value = "  +4.55%  "

$ws.Range("D31").Value = "159.11"
$ws.Range("E31").Value = "  +1.56%  "

$ws.Range("D32").Value = "33.80"
$ws.Range("E32").Value = "  +6.15%  "

$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("E34").Value = "  +3.67%  "

$ws.Range("D35").Value = "0.0741"
$ws.Range("E35").Value = "  +4.03%  "

$ws.Range("D36").Value = "3.03"
$ws.Range("E36").Value = "  +2.01%  "

$ws.Range("E37").Value = "  +2.79%  "

$ws.Range("E38").Value = "  +5.23%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "0.115"
$ws.Range("E39").Value = "  +3.62%  "

$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "16.46"
$ws.Range("E40").Value = "  +7.40%  "

$ws.Range("D41").Value = "1.78"
$ws.Range("E41").Value = "  +2.28%  "

$ws.Range("D42").Value = "3.91"
$ws.Range("E42").Value = "  +5.13%  "

$ws.Range("D43").Value = "2.049.51"
$ws.Range("E43").Value = "  -2.83%  "

$ws.Range("D44").Value = "19.61"
$ws.Range("E44").Value = "  +8.03%  "

$ws.Range("E45").Value = "  +3.28%  "

$ws.Range("E46").Value = "  +2.11%  "

$ws.Range("D47").Value = "2.06"
$ws.Range("E47").Value = "  +4.35%  "

$ws.Range("D48").Value = "2.84"
$ws.Range("E48").Value = "  +4.92%  "

$ws.Range("D49").Value = "1.51"
$ws.Range("E49").Value = "  +3.29%  "

$ws.Range("E50").Value = "  +3.59%  "

$ws.Range("D51").Value = "51.71"
$ws.Range("E51").Value = "  +6.05%  "
